$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 83 (was id 81 / match 8075530 TuRU Dsseldorf vs DV Solingen)
# now becomes match 8075296 FC Monheim vs VFB Hilden II (previously row 84's data)
$ws.Range("B83").Value = 8075296
$ws.Range("E83").Value = "FC Monheim"
$ws.Range("F83").Value = "VFB Hilden II"
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 2
$ws.Range("I83").Value = 1
$ws.Range("J83").Value = 2
$ws.Range("K83").Value = "A"
$ws.Range("L83").Value = 1.533
$ws.Range("M83").Value = 4.75
$ws.Range("N83").Value = 4
$ws.Range("O83").Value = 1.4
$ws.Range("P83").Value = 5.25
$ws.Range("Q83").Value = 5
$ws.Range("R83").Value = -1.5
$ws.Range("S83").Value = 1.975
$ws.Range("T83").Value = 1.825
$ws.Range("U83").Value = 3.75
$ws.Range("V83").Value = 1.9
$ws.Range("W83").Value = 1.9
$ws.Range("X83").Value = -1
$ws.Range("Y83").Value = -1
$ws.Range("Z83").Value = 4
$ws.Range("AA83").Value = -1
$ws.Range("AB83").Value = 0.825
$ws.Range("AC83").Value = -1
$ws.Range("AD83").Value = 0.8999999999999999

# Row 84 (was id 82 / match 8075296 FC Monheim vs VFB Hilden II)
# now becomes match 8075530 TuRU Dsseldorf vs DV Solingen (previously row 83's data)
$ws.Range("B84").Value = 8075530
$ws.Range("E84").Value = "TuRU Dsseldorf"
$ws.Range("F84").Value = "DV Solingen"
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = "H"
$ws.Range("L84").Value = 2.1
$ws.Range("M84").Value = 3.75
$ws.Range("N84").Value = 2.7
$ws.Range("O84").Value = 2.375
$ws.Range("P84").Value = 3.75
$ws.Range("Q84").Value = 2.45
$ws.Range("R84").Value = 0
$ws.Range("S84").Value = 1.85
$ws.Range("T84").Value = 1.95
$ws.Range("U84").Value = 3
$ws.Range("V84").Value = 1.85
$ws.Range("W84").Value = 1.95
$ws.Range("X84").Value = 1.375
$ws.Range("Y84").Value = -1
$ws.Range("Z84").Value = -1
$ws.Range("AA84").Value = 0.8500000000000001
$ws.Range("AB84").Value = -1
$ws.Range("AC84").Value = -1
$ws.Range("AD84").Value = 0.95

# Row 129 (was id 127 / match 8271342 BSC Hastedt vs SV Grohn)
# now becomes match 8271343 ASV Mettmann vs SC Dsseldorf West (previously row 130's data)
$ws.Range("B129").Value = 8271343
$ws.Range("E129").Value = "ASV Mettmann"
$ws.Range("F129").Value = "SC Dsseldorf West"
$ws.Range("G129").Value = 2
$ws.Range("H129").Value = 1
$ws.Range("K129").Value = "H"
$ws.Range("L129").Value = 2.25
$ws.Range("M129").Value = 3.5
$ws.Range("N129").Value = 2.6
$ws.Range("O129").Value = 2
$ws.Range("P129").Value = 3.9
$ws.Range("Q129").Value = 2.75
$ws.Range("R129").Value = -0.25
$ws.Range("S129").Value = 1.85
$ws.Range("T129").Value = 1.95
$ws.Range("U129").Value = 4.25
$ws.Range("V129").Value = 1.925
$ws.Range("W129").Value = 1.875
$ws.Range("X129").Value = 1
$ws.Range("Y129").Value = -1
$ws.Range("Z129").Value = -1
$ws.Range("AA129").Value = 0.8500000000000001
$ws.Range("AB129").Value = -1
$ws.Range("AC129").Value = -1
$ws.Range("AD129").Value = 0.875

# Row 130 (was id 128 / match 8271343 ASV Mettmann vs SC Dsseldorf West)
# now becomes match 8271342 BSC Hastedt vs SV Grohn (previously row 129's data)
$ws.Range("B130").Value = 8271342
$ws.Range("E130").Value = "BSC Hastedt"
$ws.Range("F130").Value = "SV Grohn"
$ws.Range("G130").Value = 2
$ws.Range("H130").Value = 3
$ws.Range("K130").Value = "A"
$ws.Range("L130").Value = 3.1
$ws.Range("M130").Value = 4
$ws.Range("N130").Value = 1.833
$ws.Range("O130").Value = 2.5
$ws.Range("P130").Value = 4.333
$ws.Range("Q130").Value = 2.05
$ws.Range("R130").Value = 0.25
$ws.Range("S130").Value = 1.875
$ws.Range("T130").Value = 1.925
$ws.Range("U130").Value = 4.5
$ws.Range("V130").Value = 1.85
$ws.Range("W130").Value = 1.95
$ws.Range("X130").Value = -1
$ws.Range("Y130").Value = -1
$ws.Range("Z130").Value = 1.05
$ws.Range("AA130").Value = -1
$ws.Range("AB130").Value = 0.925
$ws.Range("AC130").Value = 0.8500000000000001
$ws.Range("AD130").Value = -1
